# Fruta / hortaliza, semanal
# Adds a new weekly price observation at the top of the Pomelo (Start Ruby,
# Primera) history block (rows 32-133) by inserting a new row at position 32
# and pushing the existing rows down by one (old row 133 becomes row 134).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 32; everything from 32-133 shifts down
# to 33-134 (carrying its formatting, including the date style on column D).
$ws.Rows("32").Insert()

# Populate the newly inserted row with this week's observation.
$ws.Range("A32").Value = 10
$ws.Range("B32").Value = "Vega Modelo de Temuco"
$ws.Range("C32").Value = "La Araucanía"
$ws.Range("D32").Value = 44453
$ws.Range("E32").Value = 9
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100102
$ws.Range("H32").Value = "Cítricos"
$ws.Range("I32").Value = 100102006
$ws.Range("J32").Value = "Pomelo"
$ws.Range("K32").Value = "Start Ruby"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 55
$ws.Range("N32").Value = 13000
$ws.Range("O32").Value = 13000
$ws.Range("P32").Value = 13000
$ws.Range("Q32").Value = "$/bandeja 15 kilos granel"
$ws.Range("R32").Value = "Región de O'Higgins"
$ws.Range("S32").Value = 867
$ws.Range("T32").Value = 15
